$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Myoc"
$ws.Range("C2").Value = "Fzd3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1664173333333333
$ws.Range("H2").Value = 0.499252
$ws.Range("I2").Value = 0.009928329481286188
$ws.Range("J2").Value = 0.01026989293949489
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2078313333333333
$ws.Range("N2").Value = 0.623494
$ws.Range("O2").Value = 0.08621557350328635
$ws.Range("P2").Value = 0.112461889302165
$ws.Range("Q2").Value = 0.03458673627644444
$ws.Range("R2").Value = 0.311280626488
$ws.Range("S2").Value = 0.0008559766201586741
$ws.Range("T2").Value = 0.00115497156290656

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Myoc"
$ws.Range("C3").Value = "Fzd3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1664173333333333
$ws.Range("H3").Value = 0.499252
$ws.Range("I3").Value = 0.009928329481286188
$ws.Range("J3").Value = 0.01026989293949489
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.498127
$ws.Range("N3").Value = 1.494381
$ws.Range("O3").Value = 0.206640184103479
$ws.Range("P3").Value = 0.2695469573039334
$ws.Range("Q3").Value = 0.08289696700133334
$ws.Range("R3").Value = 0.746072703012
$ws.Range("S3").Value = 0.002051591831852976
$ws.Range("T3").Value = 0.002768218393677996

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Myoc"
$ws.Range("C4").Value = "Fzd3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1664173333333333
$ws.Range("H4").Value = 0.499252
$ws.Range("I4").Value = 0.009928329481286188
$ws.Range("J4").Value = 0.01026989293949489
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.016887
$ws.Range("N4").Value = 0.050661
$ws.Range("O4").Value = 0.007005307459654767
$ws.Range("P4").Value = 0.009137909545139137
$ws.Range("Q4").Value = 0.002810289508
$ws.Range("R4").Value = 0.025292605572
$ws.Range("S4").Value = 0.00006955100057716447
$ws.Range("T4").Value = 0.00009384535271936738

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Myoc"
$ws.Range("C5").Value = "Fzd3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1664173333333333
$ws.Range("H5").Value = 0.499252
$ws.Range("I5").Value = 0.009928329481286188
$ws.Range("J5").Value = 0.01026989293949489
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.6877555
$ws.Range("N5").Value = 3.375511
$ws.Range("O5").Value = 0.7001389349335798
$ws.Range("P5").Value = 0.6088532438487625
$ws.Range("Q5").Value = 0.2808717696286667
$ws.Range("R5").Value = 1.685230617772
$ws.Range("S5").Value = 0.006951210028697372
$ws.Range("T5").Value = 0.006252857630190966

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Myoc"
$ws.Range("C6").Value = "Fzd3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 14.92301466666667
$ws.Range("H6").Value = 44.76904399999999
$ws.Range("I6").Value = 0.8902955208876449
$ws.Range("J6").Value = 0.920924280490686
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2078313333333333
$ws.Range("N6").Value = 0.623494
$ws.Range("O6").Value = 0.08621557350328635
$ws.Range("P6").Value = 0.112461889302165
$ws.Range("Q6").Value = 3.101470035526222
$ws.Range("R6").Value = 27.91323031973599
$ws.Range("S6").Value = 0.07675733892073534
$ws.Range("T6").Value = 0.1035688844882195

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Myoc"
$ws.Range("C7").Value = "Fzd3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 14.92301466666667
$ws.Range("H7").Value = 44.76904399999999
$ws.Range("I7").Value = 0.8902955208876449
$ws.Range("J7").Value = 0.920924280490686
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.498127
$ws.Range("N7").Value = 1.494381
$ws.Range("O7").Value = 0.206640184103479
$ws.Range("P7").Value = 0.2695469573039334
$ws.Range("Q7").Value = 7.433556526862667
$ws.Range("R7").Value = 66.90200874176399
$ws.Range("S7").Value = 0.1839708303427257
$ws.Range("T7").Value = 0.2482323377135785

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Myoc"
$ws.Range("C8").Value = "Fzd3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.92301466666667
$ws.Range("H8").Value = 44.76904399999999
$ws.Range("I8").Value = 0.8902955208876449
$ws.Range("J8").Value = 0.920924280490686
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.016887
$ws.Range("N8").Value = 0.050661
$ws.Range("O8").Value = 0.007005307459654767
$ws.Range("P8").Value = 0.009137909545139137
$ws.Range("Q8").Value = 0.252004948676
$ws.Range("R8").Value = 2.268044538084
$ws.Range("S8").Value = 0.006236793853771445
$ws.Range("T8").Value = 0.00841532277304623

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Myoc"
$ws.Range("C9").Value = "Fzd3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.92301466666667
$ws.Range("H9").Value = 44.76904399999999
$ws.Range("I9").Value = 0.8902955208876449
$ws.Range("J9").Value = 0.920924280490686
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.6877555
$ws.Range("N9").Value = 3.375511
$ws.Range("O9").Value = 0.7001389349335798
$ws.Range("P9").Value = 0.6088532438487625
$ws.Range("Q9").Value = 25.18640008024733
$ws.Range("R9").Value = 151.118400481484
$ws.Range("S9").Value = 0.6233305577704124
$ws.Range("T9").Value = 0.5607077355158417

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Myoc"
$ws.Range("C10").Value = "Fzd3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.6724345
$ws.Range("H10").Value = 3.344869
$ws.Range("I10").Value = 0.09977614963106883
$ws.Range("J10").Value = 0.06880582656981911
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2078313333333333
$ws.Range("N10").Value = 0.623494
$ws.Range("O10").Value = 0.08621557350328635
$ws.Range("P10").Value = 0.112461889302165
$ws.Range("Q10").Value = 0.3475842920476667
$ws.Range("R10").Value = 2.085505752286
$ws.Range("S10").Value = 0.008602257962392312
$ws.Range("T10").Value = 0.007738033251038959

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Myoc"
$ws.Range("C11").Value = "Fzd3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.6724345
$ws.Range("H11").Value = 3.344869
$ws.Range("I11").Value = 0.09977614963106883
$ws.Range("J11").Value = 0.06880582656981911
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.498127
$ws.Range("N11").Value = 1.494381
$ws.Range("O11").Value = 0.206640184103479
$ws.Range("P11").Value = 0.2695469573039334
$ws.Range("Q11").Value = 0.8330847801815001
$ws.Range("R11").Value = 4.998508681089
$ws.Range("S11").Value = 0.02061776192890033
$ws.Range("T11").Value = 0.01854640119667688

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Myoc"
$ws.Range("C12").Value = "Fzd3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.6724345
$ws.Range("H12").Value = 3.344869
$ws.Range("I12").Value = 0.09977614963106883
$ws.Range("J12").Value = 0.06880582656981911
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.016887
$ws.Range("N12").Value = 0.050661
$ws.Range("O12").Value = 0.007005307459654767
$ws.Range("P12").Value = 0.009137909545139137
$ws.Range("Q12").Value = 0.0282424014015
$ws.Range("R12").Value = 0.169454408409
$ws.Range("S12").Value = 0.0006989626053061567
$ws.Range("T12").Value = 0.000628741419373538

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Myoc"
$ws.Range("C13").Value = "Fzd3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.6724345
$ws.Range("H13").Value = 3.344869
$ws.Range("I13").Value = 0.09977614963106883
$ws.Range("J13").Value = 0.06880582656981911
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.6877555
$ws.Range("N13").Value = 3.375511
$ws.Range("O13").Value = 0.7001389349335798
$ws.Range("P13").Value = 0.6088532438487625
$ws.Range("Q13").Value = 2.82266052576475
$ws.Range("R13").Value = 11.290642103059
$ws.Range("S13").Value = 0.06985716713447003
$ws.Range("T13").Value = 0.04189265070272973
